$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "2021年"
$ws.Range("A11").Style = $ws.Range("A10").Style

$values = @{
    "B11" = 447.57
    "C11" = 130.66
    "D11" = 75.51000000000001
    "F11" = 93.59
    "G11" = 924.96
    "H11" = 59.33
    "I11" = 352.14
    "J11" = 79.73999999999999
    "K11" = 16694.74
    "L11" = 19.88
    "M11" = 26.39
    "N11" = 88.89
    "O11" = 12.3
    "P11" = 502.05
    "Q11" = 156.83
    "R11" = 10.41
    "S11" = 84.65000000000001
    "T11" = 308.91
    "U11" = 2181.26
    "V11" = 841.99
    "W11" = 1506.21
    "X11" = 157.74
    "Y11" = 1407.27
    "Z11" = 392.83
    "AA11" = 6.7
    "AB11" = 639.3099999999999
    "AC11" = 730.91
    "AD11" = 33.23
    "AE11" = 57
    "AF11" = 1430.4
    "AG11" = 522.33
    "AH11" = 52.73
    "AI11" = 279.04
    "AJ11" = 100.43
    "AK11" = 251.67
    "AL11" = 842.9400000000001
    "AM11" = 513.41
    "AN11" = 76.63
    "AO11" = 83.98
    "AP11" = 1038.47
    "AQ11" = 173.33
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$ws.Range("E11").Value = ""
